$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 160 (pushing old rows 160-228 down to 162-230)
$ws.Rows("160:161").Insert()

# Fill in the new row 160 with its data
$ws.Cells.Item(160, 1).Value = 3
$ws.Cells.Item(160, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44489
$ws.Cells.Item(160, 5).Value = 5
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100108
$ws.Cells.Item(160, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(160, 9).Value = 100108002
$ws.Cells.Item(160, 10).Value = "Mango"
$ws.Cells.Item(160, 11).Value = "Sin especificar"
$ws.Cells.Item(160, 12).Value = "Primera"
$ws.Cells.Item(160, 13).Value = 456
$ws.Cells.Item(160, 14).Value = 7500
$ws.Cells.Item(160, 15).Value = 7500
$ws.Cells.Item(160, 16).Value = 7500
$ws.Cells.Item(160, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(160, 18).Value = "Perú"
$ws.Cells.Item(160, 19).Value = 1875
$ws.Cells.Item(160, 20).Value = 4

# Fill in the new row 161 with its data
$ws.Cells.Item(161, 1).Value = 3
$ws.Cells.Item(161, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44489
$ws.Cells.Item(161, 5).Value = 5
$ws.Cells.Item(161, 6).Value = "Fruta"
$ws.Cells.Item(161, 7).Value = 100108
$ws.Cells.Item(161, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(161, 9).Value = 100108002
$ws.Cells.Item(161, 10).Value = "Mango"
$ws.Cells.Item(161, 11).Value = "Sin especificar"
$ws.Cells.Item(161, 12).Value = "Segunda"
$ws.Cells.Item(161, 13).Value = 228
$ws.Cells.Item(161, 14).Value = 7500
$ws.Cells.Item(161, 15).Value = 7500
$ws.Cells.Item(161, 16).Value = 7500
$ws.Cells.Item(161, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(161, 18).Value = "Perú"
$ws.Cells.Item(161, 19).Value = 1875
$ws.Cells.Item(161, 20).Value = 4
